# Update "想去人数" (interested-people count) figures that changed between
# crawls, on both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 15673
$ws1.Range("F9").Value  = 15432
$ws1.Range("F11").Value = 9031
$ws1.Range("F14").Value = 1012
$ws1.Range("F15").Value = 92
$ws1.Range("F39").Value = 5556

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 15673
$ws4.Range("F9").Value  = 15432
$ws4.Range("F11").Value = 9031
$ws4.Range("F14").Value = 1012
$ws4.Range("F15").Value = 92
$ws4.Range("F41").Value = 5556
